# Add a new weekly record at the top of the Femacal de La Calera - Albahaca block.
# This shifts existing rows 133:161 down to 134:162 and fills the new row 133
# with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 133; Excel shifts rows 133:161 down to 134:162 and the
# used range grows from A1:R161 to A1:R162 automatically.
$ws.Rows.Item(133).Insert()

# Populate the newly inserted row 133 with the new weekly record.
$ws.Range("A133").Value = 3
$ws.Range("B133").Value = "Femacal de La Calera"
$ws.Range("C133").Value = "Coquimbo"
$ws.Range("D133").Value = 44637
$ws.Range("E133").Value = 5
$ws.Range("F133").Value = 100112052
$ws.Range("G133").Value = "Albahaca"
$ws.Range("H133").Value = "Sin especificar"
$ws.Range("I133").Value = "Primera"
$ws.Range("J133").Value = 150
$ws.Range("K133").Value = 4000
$ws.Range("L133").Value = 4500
$ws.Range("M133").Value = 4267
$ws.Range("N133").Value = "$/docena de matas"
$ws.Range("O133").Value = "Provincia de Quillota"
$ws.Range("P133").Value = 711
$ws.Range("Q133").Value = 6
$ws.Range("R133").Value = "Hortaliza"
